# Add a new "Skill Description" column between SkillCode (A) and SFIA Level (B).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; existing B (SFIA Level), C (Keycode), D (Description)
# shift right to C, D, E respectively.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Skill Description"

# Map each SkillCode (column A) to its human readable "Skill Description".
$skillNames = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "ETDL"       = "Learning delivery"
    "MADE"       = "MADE"
    "KNOW"       = "Knowledge management"
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -and $skillNames.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $skillNames[$code]
    }
}
